# Derek's Log.xlsx - add Friday (Oct 14, 2016 / serial 42652) entries to the
# "Logs" sheet, following the existing day-block layout/formatting used
# throughout the sheet (a blank "day divider" row followed by the day's
# task rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# ---------------------------------------------------------------------
# 1. Build out rows 322-334 by cloning existing, identically formatted
#    rows (this carries over all cell styles/borders/fills exactly) and
#    then overwriting the cell values that differ.
# ---------------------------------------------------------------------

# Row 322: blank "day divider" row (style pattern matches e.g. row 313)
$ws.Range("A313:F313").Copy($ws.Range("A322:F322"))
$ws.Range("C322").Value = "FRIDAY"

# Row 323: AV Shutdown
$ws.Range("A3:F3").Copy($ws.Range("A323:F323"))
$ws.Range("A323").Value = "AV Shutdown"
$ws.Range("B323").Value = 42652
$ws.Range("C323").Value = "1530"
$ws.Range("D323").Value = "R"
$ws.Range("E323").Value = "S203"
$ws.Range("F323").ClearContents()

# Row 324: AV Shutdown
$ws.Range("A3:F3").Copy($ws.Range("A324:F324"))
$ws.Range("A324").Value = "AV Shutdown"
$ws.Range("B324").Value = 42652
$ws.Range("C324").Value = "1530"
$ws.Range("D324").Value = "R"
$ws.Range("E324").Value = "N203"
$ws.Range("F324").ClearContents()

# Row 325: Other (style pattern, and F value "Door code 11012*", match
# row 74 exactly - only the date and time change)
$ws.Range("A74:F74").Copy($ws.Range("A325:F325"))
$ws.Range("B325").Value = 42652
$ws.Range("C325").Value = "1730"

# Row 326: Pickup PC
$ws.Range("A3:F3").Copy($ws.Range("A326:F326"))
$ws.Range("A326").Value = "Pickup PC"
$ws.Range("B326").Value = 42652
$ws.Range("C326").Value = "1700"
$ws.Range("D326").Value = "VC"
$ws.Range("E326").Value = "001-JCR"
$ws.Range("F326").Value = "Pick up roll in PC and Projector carts. Return to Vanier 040 basement storeroom. Key is in Founders 164 storeroom."
$ws.Range("A326:F326").RowHeight = 30

# Row 327: Demo
$ws.Range("A3:F3").Copy($ws.Range("A327:F327"))
$ws.Range("A327").Value = "Demo"
$ws.Range("B327").Value = 42652
$ws.Range("C327").Value = "1600"
$ws.Range("D327").Value = "FC"
$ws.Range("E327").Value = "203"
$ws.Range("F327").Value = "Demo built in PC to client. Make sure client is happy."

# Row 328: Demo
$ws.Range("A3:F3").Copy($ws.Range("A328:F328"))
$ws.Range("A328").Value = "Demo"
$ws.Range("B328").Value = 42652
$ws.Range("C328").Value = "1645"
$ws.Range("D328").Value = "MC"
$ws.Range("E328").Value = "140-SCR"
$ws.Range("F328").Value = "Door code 7083*. PC AND PROJECTOR IN ROOM ALREADY. JUST DEMO TO CLIENT."

# Row 329: AV Shutdown
$ws.Range("A3:F3").Copy($ws.Range("A329:F329"))
$ws.Range("A329").Value = "AV Shutdown"
$ws.Range("B329").Value = 42652
$ws.Range("C329").Value = "1730"
$ws.Range("D329").Value = "R"
$ws.Range("E329").Value = "N102"
$ws.Range("F329").Value = "Nat Taylor Cinema. Lock cinema all doors after shutdown."

# Row 330: Other (style pattern, and F value "Door code 11012*", match
# row 74 exactly - only the date and time change)
$ws.Range("A74:F74").Copy($ws.Range("A330:F330"))
$ws.Range("B330").Value = 42652
$ws.Range("C330").Value = "2030"

# Row 331: Pickup PC
$ws.Range("A3:F3").Copy($ws.Range("A331:F331"))
$ws.Range("A331").Value = "Pickup PC"
$ws.Range("B331").Value = 42652
$ws.Range("C331").Value = "1930"
$ws.Range("D331").Value = "MC"
$ws.Range("E331").Value = "140-SCR"
$ws.Range("F331").Value = "Door code 7083*. Pick up PC and Projector from Mac SCR. Pick up portable screen and return all equipment to Fouders 156A storeroom. "
$ws.Range("A331:F331").RowHeight = 45

# Row 332: AV Shutdown
$ws.Range("A3:F3").Copy($ws.Range("A332:F332"))
$ws.Range("A332").Value = "AV Shutdown"
$ws.Range("B332").Value = 42652
$ws.Range("C332").Value = "2030"
$ws.Range("D332").Value = "R"
$ws.Range("E332").Value = "S203"
$ws.Range("F332").ClearContents()

# Row 333: AV Shutdown
$ws.Range("A3:F3").Copy($ws.Range("A333:F333"))
$ws.Range("A333").Value = "AV Shutdown"
$ws.Range("B333").Value = 42652
$ws.Range("C333").Value = "1730"
$ws.Range("D333").Value = "R"
$ws.Range("E333").Value = "N203"
$ws.Range("F333").ClearContents()

# Row 334: AV Shutdown
$ws.Range("A3:F3").Copy($ws.Range("A334:F334"))
$ws.Range("A334").Value = "AV Shutdown"
$ws.Range("B334").Value = 42652
$ws.Range("C334").Value = "1900"
$ws.Range("D334").Value = "FC"
$ws.Range("E334").Value = "203"
$ws.Range("F334").ClearContents()

# ---------------------------------------------------------------------
# 2. Update the sheet view: scroll the frozen pane down and move the
#    selection to the newly added rows.
# ---------------------------------------------------------------------
$ws.Activate()
$w = $excel.ActiveWindow
$w.ScrollRow = 315
$ws.Range("B324:B334").Select()
